$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reassign the "Periodo Mora" (column E) values for rows 16-23 so the
# table now lists periods 2207-2212 in order, alternating the two
# workers for the shared periods 2207 and 2208.
$ws.Range("C16").Value = "40987503"
$ws.Range("D16").Value = "GLORIA VICTORIA VANEGAS REYES"
$ws.Range("E16").Value = "2207"

$ws.Range("C17").Value = "1043964778"
$ws.Range("D17").Value = "NELSON ENRIQUE PACHECO BOHORQUEZ"
$ws.Range("E17").Value = "2207"

$ws.Range("C18").Value = "40987503"
$ws.Range("D18").Value = "GLORIA VICTORIA VANEGAS REYES"
$ws.Range("E18").Value = "2208"

$ws.Range("C19").Value = "1043964778"
$ws.Range("D19").Value = "NELSON ENRIQUE PACHECO BOHORQUEZ"
$ws.Range("E19").Value = "2208"

$ws.Range("C20").Value = "40987503"
$ws.Range("D20").Value = "GLORIA VICTORIA VANEGAS REYES"
$ws.Range("E20").Value = "2209"

$ws.Range("C21").Value = "40987503"
$ws.Range("D21").Value = "GLORIA VICTORIA VANEGAS REYES"
$ws.Range("E21").Value = "2210"

$ws.Range("C22").Value = "40987503"
$ws.Range("D22").Value = "GLORIA VICTORIA VANEGAS REYES"
$ws.Range("E22").Value = "2211"

$ws.Range("C23").Value = "40987503"
$ws.Range("D23").Value = "GLORIA VICTORIA VANEGAS REYES"
$ws.Range("E23").Value = "2212"
